# Weekly fruit/vegetable price update:
# A new weekly record is inserted before the existing row 338, pushing the
# remaining records (old rows 338-362) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(338).Insert()

$ws.Cells.Item(338, 1).Value = 4
$ws.Cells.Item(338, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(338, 3).Value = "Los Lagos"
$ws.Cells.Item(338, 4).Value = 44826
$ws.Cells.Item(338, 5).Value = 10
$ws.Cells.Item(338, 6).Value = 100112045
$ws.Cells.Item(338, 7).Value = "Zapallo"
$ws.Cells.Item(338, 8).Value = "Paine"
$ws.Cells.Item(338, 9).Value = "1a (guarda)"
$ws.Cells.Item(338, 10).Value = 500
$ws.Cells.Item(338, 11).Value = 600
$ws.Cells.Item(338, 12).Value = 600
$ws.Cells.Item(338, 13).Value = 600
$ws.Cells.Item(338, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(338, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(338, 16).Value = 600
$ws.Cells.Item(338, 17).Value = 1
$ws.Cells.Item(338, 18).Value = "Hortaliza"
